$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44819
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25500
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("S2").Value = 1417

# Row 3
$ws.Range("D3").Value = 44316
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 1111

# Row 5
$ws.Range("D5").Value = 44687
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 19000
$ws.Range("P5").Value = 18500
$ws.Range("S5").Value = 1028

# Row 6
$ws.Range("D6").Value = 45044
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("S6").Value = 972

# Row 7
$ws.Range("D7").Value = 45030
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15500
$ws.Range("Q7").Value = "`$/caja 18 kilos granel"
$ws.Range("S7").Value = 861

# Row 8
$ws.Range("D8").Value = 44516
$ws.Range("N8").Value = 33000
$ws.Range("O8").Value = 34000
$ws.Range("P8").Value = 33500
$ws.Range("S8").Value = 1861

# Row 9
$ws.Range("D9").Value = 45002
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 12500
$ws.Range("Q9").Value = "`$/caja 18 kilos"
$ws.Range("S9").Value = 694

# Row 10
$ws.Range("D10").Value = 45014
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 13600
$ws.Range("Q10").Value = "`$/caja 18 kilos"
$ws.Range("S10").Value = 756

# Row 11
$ws.Range("D11").Value = 45014
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("S11").Value = 556

# Row 12
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14500
$ws.Range("S12").Value = 806

# Row 13
$ws.Range("D13").Value = 44280
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 667

# Row 14
$ws.Range("D14").Value = 44699
$ws.Range("M14").Value = 100
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21000
$ws.Range("S14").Value = 1167

# Row 15
$ws.Range("D15").Value = 44699
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 50
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 18000
$ws.Range("S15").Value = 1000

